$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" date fields: 22/03/2017 -> 23/03/2017
#    (these live on the slide master, the 11 slide layouts and the notes
#    master - there are no such fields on the slides themselves).
# ---------------------------------------------------------------------------

function Update-DateShape($shape) {
    if (-not $shape.HasTextFrame) { return }
    $tr = $shape.TextFrame.TextRange
    $t = $tr.Text
    if ($t -eq "3/22/2017") {
        $tr.Text = "3/23/2017"
    } elseif ($t -eq "22/03/2017") {
        $tr.Text = "23/03/2017"
    }
}

# Slide master
$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    Update-DateShape $master.Shapes.Item($j)
}

# Every slide layout attached to the master
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        Update-DateShape $layout.Shapes.Item($j)
    }
}

# Notes master (HasNotesMaster does not reliably surface as a PowerShell
# boolean in this host, so just grab NotesMaster directly).
$notesMaster = $p.NotesMaster
for ($j = 1; $j -le $notesMaster.Shapes.Count; $j++) {
    Update-DateShape $notesMaster.Shapes.Item($j)
}

# ---------------------------------------------------------------------------
# 2) Slide 1 title: split the run so a new phrase "Protótipo de " is inserted
#    right before "Elevador".
#       "Automação de um Elevador Residencial..."
#    -> "Automação de um Protótipo de Elevador Residencial..."
# ---------------------------------------------------------------------------

$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

$oldWord = "Elevador "
$newPhrase = "Protótipo de Elevador "

$fullText = $titleRange.Text
$pos = $fullText.IndexOf($oldWord)
if ($pos -ge 0) {
    # TextRange.Characters is 1-based; .Text is a plain 0-based .NET string
    # over the very same character sequence (including the leading line
    # break), so Characters-index = .Text-0-based-index + 1.
    $startIndex = $pos + 1
    $target = $titleRange.Characters($startIndex, $oldWord.Length)
    $target.Text = $newPhrase
}
